$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0,1.186047627950188,0.2387000123951566,0.2226785005468557,0,1.192283006441741,0.6132414258254641,0.7327866179132911,0,0.1981633092984048,0,0.3616956605820718,0,1.246523001320917,2.676150110192594),
    @(1,1.093298377139433,0.2319060946352209,0.2219172207628901,0,1.194552975616624,0.6132278077361661,0.7367124965767431,0,0.1999377577772492,0,0.3538042235259695,0,1.252480024247184,2.684057865062783),
    @(2,1.036495716730258,0.2277245149831799,0.2215245335469049,0,1.19659359528567,0.6136372241045507,0.7394537970780988,0,0.2011021738988532,0,0.349085050130995,0,1.256596613984257,2.690481307476304),
    @(3,1.013386627025113,0.2260181159400076,0.2213833737355628,0,1.197587911280593,0.6139090671336476,0.7406541449356254,0,0.2015955312009279,0,0.3471938874005644,0,1.258389766659207,2.693493230785222),
    @(4,1.009551756048239,0.2257346313128181,0.2213610755903659,0,1.197762848896822,0.613960546857335,0.7408584915164909,0,0.2016785916266386,0,0.3468817963103703,0,1.258694506788949,2.694017174812316),
    @(5,1.036183901367878,0.2277015112716754,0.2215225533472491,0,1.196606345902111,0.6136404651916649,0.7394696482453824,0,0.2011087511524536,0,0.3490594156585729,0,1.256620328702681,2.690520330723189),
    @(6,1.154038487284197,0.2363596598053306,0.2224005394045747,0,1.192931490516365,0.6131499879294893,0.734071642457792,0,0.1987596030585514,0,0.3589486196641616,0,1.248481845742589,2.678551262043214),
    @(7,1.386239741241013,0.2532517840055277,0.2247125371414427,0,1.190854808565767,0.6155062941612357,0.7261083596243196,0,0.194746528007439,0,0.3793349204256202,0,1.236155619766777,2.667523649625991),
    @(8,1.557427720800376,0.2656022257149289,0.2267676642471628,0,1.192454414898435,0.6192661537461817,0.7218533694610727,0,0.192159024044269,0,0.3949098777492566,0,1.229304213377944,2.667014416375935),
    @(9,1.635418960461777,0.2712061792497593,0.2277793219759587,0,1.193860484398982,0.6214184548965989,0.7202635614067674,0,0.191060062480588,0,0.4021232841330118,0,1.226664008957499,2.668433157287751),
    @(10,1.664967422387633,0.2733260387963696,0.2281733853198205,0,1.194490422172152,0.6222971034433442,0.7197112162682942,0,0.1906551317880911,0,0.4048730769649183,0,1.225732589734093,2.669207788392868),
    @(11,1.658603001965105,0.272869591360319,0.2280880299080366,0,1.194350418939109,0.6221050401786954,0.7198279648495287,0,0.1907418418977969,0,0.4042800523972119,0,1.225930149244505,2.66903039904895),
    @(12,1.637849643295738,0.2713806271913484,0.2278115223339086,0,1.193910356285357,0.6214894664236681,0.7202171242286965,0,0.1910265237330933,0,0.4023491470667153,0,1.226586011201888,2.66849212877753),
    @(13,1.625139490895322,0.2704682972275521,0.227643580100704,0,1.193653499349381,0.62112069685422,0.7204619641065193,0,0.1912023607594655,0,0.4011687798381303,0,1.226996644675651,2.668193338926329),
    @(14,1.552332997003987,0.2652356908516538,0.2267030889535846,0,1.192376170707476,0.6191343945030212,0.7219642217585829,0,0.192232414665952,0,0.3944410276949952,0,1.229486332586745,2.666954913664199),
    @(15,1.507697142864345,0.2620218619301795,0.2261457416964703,0,1.191766280834983,0.6180291001211913,0.7229743451282076,0,0.1928843176649941,0,0.3903464729007737,0,1.231135614879364,2.666617928957663),
    @(16,1.48203490606528,0.2601720147549145,0.225832400379872,0,1.191479346437347,0.6174349556750087,0.7235878932412447,0,0.1932666273878034,0,0.3880034804161312,0,1.232129101184924,2.666579443706894),
    @(17,1.473348083919404,0.2595454639009205,0.225727552268097,0,1.191393164991347,0.6172409297973047,0.7238012223134405,0,0.1933973339563941,0,0.3872122659220167,0,1.232473188633421,2.666593091087549),
    @(18,1.512447568363143,0.2623641192397201,0.2262043244765408,0,1.191824596032276,0.6181424554971215,0.7228634471084376,0,0.192814160605451,0,0.3907810958277196,0,1.23095540398225,2.666637724059427),
    @(19,1.64394501969133,0.2718180341483958,0.2278924421232063,0,1.194036968020342,0.6216685482024928,0.7201014707339226,0,0.1909426012853892,0,0.4029158076526045,0,1.226391514357488,2.668643788786085),
    @(20,1.729971989852402,0.2779835994066673,0.2290596372161389,0,1.196051095813601,0.6243439047073025,0.7185859150801264,0,0.1897848363113255,0,0.4109527162834041,0,1.223807185508264,2.671338553702242),
    @(21,1.684050570441002,0.2746941792368034,0.2284308579783811,0,1.194924147436126,0.6228820606029188,0.7193683165391462,0,0.1903967756128324,0,0.4066536227943942,0,1.225150083924589,2.669773680499105),
    @(22,1.510299902370377,0.2622093914687014,0.2261778171277058,0,1.191798033308004,0.6180910788779954,0.7229134818832534,0,0.1928458551813073,0,0.3905845684561626,0,1.231036736235154,2.666628291099357),
    @(23,1.323314254206025,0.2486920509061861,0.2240242338011527,0,1.190867625631583,0.6145130122811224,0.7279822179216353,0,0.1957687231768919,0,0.3793349204256202,0,1.239102301335336,2.669174022939785)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
